# Insert a new weekly data row for "Feria Lagunitas de Puerto Montt - Cebolla"
# at row 631, shifting the existing rows 631-672 down to 632-673.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(631).Insert()

$ws.Cells.Item(631, 1).Value = 4
$ws.Cells.Item(631, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(631, 3).Value = "Los Lagos"
$ws.Cells.Item(631, 4).Value = 44826
$ws.Cells.Item(631, 5).Value = 10
$ws.Cells.Item(631, 6).Value = 100112004
$ws.Cells.Item(631, 7).Value = "Cebolla"
$ws.Cells.Item(631, 8).Value = "Sin especificar"
$ws.Cells.Item(631, 9).Value = "1a (guarda)"
$ws.Cells.Item(631, 10).Value = 750
$ws.Cells.Item(631, 11).Value = 14000
$ws.Cells.Item(631, 12).Value = 14000
$ws.Cells.Item(631, 13).Value = 14000
$ws.Cells.Item(631, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(631, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(631, 16).Value = 778
$ws.Cells.Item(631, 17).Value = 18
$ws.Cells.Item(631, 18).Value = "Hortaliza"
